$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format before writing numeric-looking
# strings, so Excel does not silently convert them to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '271.52'

# Row 3
$ws.Range("D3").Value = '22.80'

# Row 4
$ws.Range("D4").Value = '6.343'

# Row 5
$ws.Range("D5").Value = '0.06198'

# Row 6
$ws.Range("D6").Value = '3.649'

# Row 7
$ws.Range("D7").Value = '6.695'

# Row 8
$ws.Range("D8").Value = '1.385'

# Row 9
$ws.Range("D9").Value = '0.8307'

# Row 10
$ws.Range("D10").Value = '0.01377'

# Row 11
$ws.Range("D11").Value = '0.1600'

# Row 12
$ws.Range("D12").Value = '0.08295'

# Row 13
$ws.Range("D13").Value = '0.03448'

# Row 14
$ws.Range("D14").Value = '0.03178'

# Row 15
$ws.Range("B15").Value = 'ProBitToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D15").Value = '0.1240'
$ws.Range("E15").Value = '14ProBitTokenPROB'

# Row 16
$ws.Range("B16").Value = 'BitMartToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D16").Value = '0.09332'
$ws.Range("E16").Value = '15BitMartTokenBMX'

# Row 17
$ws.Range("B17").Value = 'MCDex'
$ws.Range("C17").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D17").Value = '3.844'
$ws.Range("E17").Value = '16MCDexMCB'

# Row 18
$ws.Range("B18").Value = 'BitForexToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D18").Value = '0.001666'
$ws.Range("E18").Value = '17BitForexTokenBF'

# Row 19
$ws.Range("B19").Value = 'CoinExToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D19").Value = '0.04741'
$ws.Range("E19").Value = '18CoinExTokenCET'

# Row 20
$ws.Range("B20").Value = 'TigerCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D20").Value = '0.006333'
$ws.Range("E20").Value = '19TigerCashTCH'

# Row 21
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").Value = '0.005665'
$ws.Range("E21").Value = '20HotbitTokenHTB'

# Row 22
$ws.Range("B22").Value = 'BitKan'
$ws.Range("C22").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D22").Value = '0.001077'
$ws.Range("E22").Value = '21BitKanKAN'

# Row 23
$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").Value = '0.0001501'
$ws.Range("E23").Value = '22NitroExNTX'

# Row 24
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").Value = '3.720'
$ws.Range("E24").Value = '23LEOLEO'

# Row 25
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").Value = '2.325'
$ws.Range("E25").Value = '24BTSETokenBTSE'

# Row 26
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").Value = '0.3348'
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'

# Row 40
$ws.Range("D40").Value = '0.04691'

# Row 41
$ws.Range("D41").Value = '0.007031'

# Row 42
$ws.Range("D42").Value = '0.1161'

# Row 43
$ws.Range("D43").Value = '0.003292'
$ws.Range("E43").Value = '42CEJICEJI'

# Row 45
$ws.Range("D45").Value = '0.00006273'

# Row 48
$ws.Range("D48").Value = '0.9204'

# Row 49
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").Value = '0.002117'
$ws.Range("E49").Value = '48BOLOBOLO'

# Row 50
$ws.Range("B50").Value = 'CryptobidCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc'
$ws.Range("D50").Value = '0.00001401'
$ws.Range("E50").Value = '49CryptobidCoinCBCWorstin24h'

# Row 51
$ws.Range("D51").Value = '0.01241'

# Restore the default (Normal) style on the Price column so no stray
# number-format styling is left behind, matching the original workbook.
$ws.Range("D2:D51").Style = "Normal"